# Update "想去人数" (want-to-go count) values in column F across sheets
# 展览 (Exhibitions), 演出 (Performances), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 252
$ws1.Range("F5").Value  = 1791
$ws1.Range("F8").Value  = 519
$ws1.Range("F9").Value  = 4753
$ws1.Range("F11").Value = 451
$ws1.Range("F13").Value = 989
$ws1.Range("F14").Value = 1289
$ws1.Range("F17").Value = 2991
$ws1.Range("F18").Value = 1820
$ws1.Range("F22").Value = 50
$ws1.Range("F24").Value = 938
$ws1.Range("F26").Value = 33
$ws1.Range("F27").Value = 2760
$ws1.Range("F28").Value = 1023
$ws1.Range("F29").Value = 2483
$ws1.Range("F30").Value = 254
$ws1.Range("F31").Value = 1339
$ws1.Range("F32").Value = 3619
$ws1.Range("F34").Value = 891
$ws1.Range("F35").Value = 433
$ws1.Range("F36").Value = 1136
$ws1.Range("F37").Value = 930
$ws1.Range("F38").Value = 1194
$ws1.Range("F40").Value = 878
$ws1.Range("F41").Value = 545
$ws1.Range("F42").Value = 172
$ws1.Range("F43").Value = 367
$ws1.Range("F44").Value = 288
$ws1.Range("F45").Value = 3495

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 888
$ws2.Range("F27").Value = 43

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 252
$ws4.Range("F6").Value  = 1791
$ws4.Range("F9").Value  = 519
$ws4.Range("F10").Value = 4753
$ws4.Range("F14").Value = 1289
$ws4.Range("F15").Value = 2991
$ws4.Range("F17").Value = 1820
$ws4.Range("F21").Value = 888
$ws4.Range("F24").Value = 50
$ws4.Range("F26").Value = 938
$ws4.Range("F28").Value = 2760
$ws4.Range("F31").Value = 1023
$ws4.Range("F32").Value = 2483
$ws4.Range("F33").Value = 1339
$ws4.Range("F34").Value = 3619
$ws4.Range("F37").Value = 891
$ws4.Range("F38").Value = 1136
$ws4.Range("F39").Value = 930
$ws4.Range("F41").Value = 1194
$ws4.Range("F42").Value = 878
$ws4.Range("F43").Value = 545
$ws4.Range("F44").Value = 367
$ws4.Range("F48").Value = 288
$ws4.Range("F49").Value = 3495
